# Apply the email-address updates described in the commit
# "fim do desafio outlook python minha resolucao"
#
# Column B (e-mail) text is updated for several rows while the underlying
# mailto: hyperlinks stay pointed at jharbes@hotmail.com (unchanged by
# this edit). Also move the active selection from B8 to B9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = "jharbes@icloud.com"
$ws.Range("B4").Value = "jorge.harbes@technipfmc.com"
$ws.Range("B5").Value = "jorgenamiharbes@gmail.com"
$ws.Range("B7").Value = "jharbes@icloud.com"
$ws.Range("B8").Value = "jorge.harbes@technipfmc.com"

$ws.Range("B9").Select()
